$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("enum")

# Insert a new row at position 59, pushing existing rows 59..100 down to 60..101
$ws.Rows.Item(59).Insert()

# Fill in the new row 59 with the "dot_amor" effect_type entry
$ws.Cells.Item(59, 1).Value = "effect_type"
$ws.Cells.Item(59, 2).Value = "dot_amor"
$ws.Cells.Item(59, 3).Formula = "=UPPER(A59)&""_""&UPPER(B59)"
$ws.Cells.Item(59, 4).Value = 29
$ws.Cells.Item(59, 5).Value = "지속 방어도 획득"
